# Updated cryptos list on Wed Oct 23 21:46:27 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.420.36"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "2.506.09"
$ws.Range("E3").Value = "  -4.83%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.99"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.96"
$ws.Range("E6").Value = "  +2.78%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D9").Value = "2.505.24"
$ws.Range("E9").Value = "  -4.84%  "
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("E12").Value = "  -3.32%  "
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D15").Value = "2.969.88"
$ws.Range("E15").Value = "  -4.58%  "
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("D17").Value = "66.313.11"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "2.510.32"
$ws.Range("E18").Value = "  -4.40%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.23"
$ws.Range("E19").Value = "  -6.19%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.67"
$ws.Range("E20").Value = "  -4.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "347.45"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.19"
$ws.Range("E22").Value = "  -2.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.62"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.59"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.00"
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D30").Value = "0.0₃0977"
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "528.47"
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("E35").Value = "  -3.99%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.45"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.60"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("E41").Value = "  -3.22%  "
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.09"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").Value = "  +3.67%  "
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.20"
$ws.Range("E47").Value = "  -3.02%  "
$ws.Range("E48").Value = "  -3.97%  "
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0273"
$ws.Range("E50").Value = "  -9.54%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.72"
$ws.Range("E51").Value = "  +1.15%  "
